# The species-occurrence records that were sitting in rows 9,10,11,12,13,14,
# 16,17,18 got re-matched against their correct source rows (row 15 is left
# untouched). Each physical row keeps its place on the sheet, but the full
# record (every populated column, A through AY) that lives in it changes.
#
# Strategy: snapshot every full row's current values first (so later writes
# never clobber data we still need to read), then write the snapshots back
# into their destination rows per the mapping below.
#
#   old row -> new row
#   9  -> 11
#   10 -> 14
#   11 -> 9
#   12 -> 17
#   13 -> 18
#   14 -> 16
#   16 -> 10
#   17 -> 12
#   18 -> 13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowValues($rowNum) {
    return $ws.Range("A" + $rowNum + ":AY" + $rowNum).Value()
}

# Snapshot the "before" contents of every row that moves.
$row9  = Get-RowValues 9
$row10 = Get-RowValues 10
$row11 = Get-RowValues 11
$row12 = Get-RowValues 12
$row13 = Get-RowValues 13
$row14 = Get-RowValues 14
$row16 = Get-RowValues 16
$row17 = Get-RowValues 17
$row18 = Get-RowValues 18

function Set-RowValues($rowNum, $values) {
    # Columns Y, Z, AA, AB hold plain text dates/times (e.g. "2023-08-16",
    # "00:00"). A bare Range.Value assignment lets Excel's COM layer
    # auto-detect and silently convert those text strings into real date
    # serial numbers. Pre-formatting that sub-range as Text *before*
    # assigning the row's values keeps them literal strings, matching the
    # source workbook's inlineStr cells.
    $ws.Range("Y" + $rowNum + ":AB" + $rowNum).NumberFormat = "@"
    $ws.Range("A" + $rowNum + ":AY" + $rowNum).Value = $values
}

Set-RowValues 11 $row9
Set-RowValues 14 $row10
Set-RowValues 9  $row11
Set-RowValues 17 $row12
Set-RowValues 18 $row13
Set-RowValues 16 $row14
Set-RowValues 10 $row16
Set-RowValues 12 $row17
Set-RowValues 13 $row18
